$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.201.30'
$ws.Range("E2").Value = '  +1.02%  '
$ws.Range("D3").Value = '1.859.48'
$ws.Range("E3").Value = '  +1.13%  '
$ws.Range("E4").Value = '  +0.65%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '239.11'
$c.ClearFormats()
$ws.Range("E5").Value = '  +3.44%  '
$ws.Range("E6").Value = '  +0.44%  '
$ws.Range("E7").Value = '  +0.61%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '42.35'
$c.ClearFormats()
$ws.Range("E8").Value = '  +6.65%  '
$ws.Range("E9").Value = '  +0.30%  '
$ws.Range("E10").Value = '  +1.36%  '
$ws.Range("E11").Value = '  +0.38%  '
$ws.Range("D12").Value = '2.128.17'
$ws.Range("E12").Value = '  +1.16%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '11.49'
$c.ClearFormats()
$ws.Range("E13").Value = '  +1.25%  '
$ws.Range("D14").Value = '1.844.27'
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("E15").Value = '  +0.61%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '4.71'
$c.ClearFormats()
$ws.Range("E16").Value = '  +1.60%  '
$ws.Range("D17").Value = '35.146.92'
$ws.Range("E17").Value = '  +0.83%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '69.89'
$c.ClearFormats()
$ws.Range("E18").Value = '  +0.30%  '
$ws.Range("D19").Value = '0.0₃0794'
$ws.Range("E19").Value = '  +1.05%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '241.03'
$c.ClearFormats()
$ws.Range("E20").Value = '  +0.17%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '12.22'
$c.ClearFormats()
$ws.Range("E21").Value = '  +0.41%  '
$ws.Range("E22").Value = '  +1.14%  '
$ws.Range("E23").Value = '  +0.46%  '
$ws.Range("E24").Value = '  +0.55%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '169.46'
$c.ClearFormats()
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '1.89'
$c.ClearFormats()
$ws.Range("E26").Value = '  +25.01%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '8.05'
$c.ClearFormats()
$ws.Range("E27").Value = '  +3.64%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '17.69'
$c.ClearFormats()
$ws.Range("E28").Value = '  +1.70%  '
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("B30").Value = 'BinanceUSD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.01'
$c.ClearFormats()
$ws.Range("E30").Value = '  +0.67%  '
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.0560'
$c.ClearFormats()
$ws.Range("E31").Value = '  +1.42%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '4.01'
$c.ClearFormats()
$ws.Range("E32").Value = '  +1.89%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '1.83'
$c.ClearFormats()
$ws.Range("E33").Value = '  +27.80%  '
$ws.Range("E34").Value = '  +2.11%  '
$ws.Range("E35").Value = '  +10.54%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.814'
$c.ClearFormats()
$ws.Range("E36").Value = '  +17.41%  '
$ws.Range("E37").Value = '  +7.85%  '
$ws.Range("E38").Value = '  +4.36%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.0201'
$c.ClearFormats()
$ws.Range("E39").Value = '  +4.25%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '89.99'
$c.ClearFormats()
$ws.Range("E40").Value = '  -1.73%  '
$ws.Range("D41").Value = '1.351.24'
$ws.Range("E41").Value = '  +0.84%  '
$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '15.03'
$c.ClearFormats()
$ws.Range("E42").Value = '  +3.22%  '
$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.0591'
$c.ClearFormats()
$ws.Range("E43").Value = '  +13.29%  '
$ws.Range("E44").Value = '  +2.22%  '
$ws.Range("E45").Value = '  +0.45%  '
$ws.Range("B46").Value = 'MXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '2.74'
$c.ClearFormats()
$ws.Range("E46").Value = '  -0.85%  '
$ws.Range("B47").Value = 'Gas'
$ws.Range("C47").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '12.27'
$c.ClearFormats()
$ws.Range("E47").Value = '  +44.00%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '6.57'
$c.ClearFormats()
$ws.Range("E48").Value = '  +4.94%  '
$ws.Range("D49").Value = '2.042.50'
$ws.Range("E49").Value = '  +1.36%  '
$ws.Range("E50").Value = '  -0.19%  '
$ws.Range("E51").Value = '  +0.66%  '
